# "Generate Report for Handback"
#
# For each localized-language sheet (zh-cn, de-de), the handback step now
# fills in the "Latest Target File" (col I) and "Latest Handback File"
# (col J) columns for every row, stamps the "Latest Handback DateTime"
# (col K), and turns the Status column (col C) + the Overview summary
# columns (E/F) from "In Translation" into "Handed back: in sync with
# en-US". The new Target-File cells get their own hyperlinks, just like
# the existing Source-File-Name column A links.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9e0791fb02e589c1b367e455f9a7ae157d784bd8/e2e/"

$file1 = "8af45f33-bc5b-423b-8790-28856c4fa88a"
$file2 = "cfecff16-4b8d-4140-997a-9a79e3457f78"

$file1Md = "$file1.md"
$file2Md = "$file2.md"

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: the per-language status columns (E, F) reflect the same
# "now handed back" status as the per-language sheets.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("J2").Value = "$file1.b9434d020143b78a88b91fee8a592dde901e0d85.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-01 02:23:34"

$wsZh.Range("J3").Value = "$file2.f06b108da99ff1ee4aea529dce1b17bc55676383.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-01 02:23:34"

# Recreate the hyperlinks in display order A2, I2, A3, I3 so the new
# "Latest Target File" links land right after the matching source-file
# link, each with its own relationship id.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$repoBase$file1Md", [Type]::Missing, [Type]::Missing, $file1Md)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$repoBase$file1Md", [Type]::Missing, [Type]::Missing, $file1Md)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$repoBase$file2Md", [Type]::Missing, [Type]::Missing, $file2Md)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$repoBase$file2Md", [Type]::Missing, [Type]::Missing, $file2Md)

$wsZh.Columns.Item(3).ColumnWidth = 29.1
$wsZh.Columns.Item(9).ColumnWidth = 39.1
$wsZh.Columns.Item(10).ColumnWidth = 39.1

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("J2").Value = "$file1.b9434d020143b78a88b91fee8a592dde901e0d85.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-01 02:23:53"

$wsDe.Range("J3").Value = "$file2.f06b108da99ff1ee4aea529dce1b17bc55676383.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-01 02:23:53"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$repoBase$file1Md", [Type]::Missing, [Type]::Missing, $file1Md)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$repoBase$file1Md", [Type]::Missing, [Type]::Missing, $file1Md)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$repoBase$file2Md", [Type]::Missing, [Type]::Missing, $file2Md)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$repoBase$file2Md", [Type]::Missing, [Type]::Missing, $file2Md)

$wsDe.Columns.Item(3).ColumnWidth = 29.1
$wsDe.Columns.Item(9).ColumnWidth = 39.1
$wsDe.Columns.Item(10).ColumnWidth = 39.1

Write-Host "Handback report generated."
